$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables predictoras")

# New column H: "long" header followed by interaction variable names
$ws.Range("H1").Value = "long"
$ws.Range("H2").Value = "m0_sexo_fac"
$ws.Range("H3").Value = "m02_fac"
$ws.Range("H4").Value = "m01_fac"
$ws.Range("H5").Value = "m0_edad"
$ws.Range("H8").Value = "s12_l"
$ws.Range("H7").Value = "s04_l"
$ws.Range("H6").Value = "s_imc_l"

# Match selection state from the diff
$ws.Range("H5").Select()
